# Refactor the "kategori" / "total_qty" columns out of the import-listing
# header row, replacing them with a single "plan_stock" column while
# leaving "customer" where it is.
#
# Before: A=Inv_id  B=part_name  C=part_number  D=customer  E=kategori  F=total_qty  G,H=(blank)
# After:  A=Inv_id  B=part_name  C=part_number  D=customer  E=plan_stock              F,G=(blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "kategori" column entirely; this shifts total_qty/blank cells
# one column to the left (F->E, G->F, H->G) and shrinks the header's
# bordered box from A:H down to A:G.
$ws.Columns("E").Delete() | Out-Null

# The old "total_qty" header (now sitting in E1) becomes "plan_stock".
$ws.Range("E1").Value = "plan_stock"

# That cell used to be the right-hand edge of the header box (kept its
# left/right borders from the old F/G boundary); the box now ends one
# column sooner, so the vertical divider in row 2 underneath it goes away.
$ws.Range("E2").Borders.Item(7).LineStyle = -4142  # xlEdgeLeft  = xlNone
$ws.Range("E2").Borders.Item(10).LineStyle = -4142 # xlEdgeRight = xlNone

# Move the active selection as recorded in the saved workbook.
$ws.Range("B6").Select() | Out-Null
